# Update the marksheet "Corr/total marks" values on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" - Right (correct) count: 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" - Right (correct) total: 21 -> 35
$ws.Range("B12").Value = 35

# Row 12 "Total" - Max column text "corr/total": "10/84" -> "35/140"
$ws.Range("E12").Value = "35/140"
